$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before P. This shifts the existing P,Q,S columns to
# Q,R,T respectively (and their shared formulas / ranges move with them).
$ws.Range("P1").EntireColumn.Insert()

# --- Value tweaks for rows 12/13 (Q5/Q6 interval boundary nudged) ---
$ws.Range("G12").Value = 11.5
$ws.Range("F13").Value = 11.5

# --- New column P header + content ---
$ws.Range("P1").Value = "starting_value3"

$ws.Range("P3").Formula = "=E3-((E3-D3)/1.7)"
$ws.Range("P4:P13").Formula = "=E4-((E4-D4)/1.7)"

$ws.Range("P14").Value = 19
$ws.Range("P15").Value = 0.5
$ws.Range("P16").Value = 0.25
$ws.Range("P17").Value = 0.75
$ws.Range("P18").Value = 1.9
$ws.Range("P19").Value = 5.5

# --- New column S content (summary strings referencing the new column P) ---
$ws.Range("S3").Formula = "=`$A3&"" = ""&ROUND(P3, 3)&"", """
$ws.Range("S4:S19").Formula = "=`$A4&"" = ""&ROUND(P4, 3)&"", """

# --- Row 22 summary formulas ---
$ws.Range("R22").Formula = "=CONCATENATE(""sigma = 12, binge = 5, "", R3,R4, R5, R6, R7, R8, R9, R10, R11, R12, R13, R14, R15, R16, R17, R18, R19)"
$ws.Range("S22").Formula = "=CONCATENATE(""sigma = 8, binge = 7, "", S3,S4, S5, S6, S7, S8, S9, S10, S11, S12, S13, S14, S15, S16, S17, S18, S19)"

# --- Selection matches the author's final cursor position ---
$ws.Range("S24").Select()
